# Fruta / hortaliza, semanal
# Rotates the weekly price-report rows: each row (2-17, except row 12 which
# stays put) takes on the D,L,M,N,O,P,Q,R,S,T values that another row held
# before the edit. Columns A,B,C,E-K (market/product metadata) are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get rotated between rows.
$cols = @("D","L","M","N","O","P","Q","R","S","T")

# target row -> source row (both refer to the *original* sheet state)
$mapping = @{
    2  = 13
    3  = 15
    4  = 6
    5  = 7
    6  = 11
    7  = 16
    8  = 2
    9  = 3
    10 = 14
    11 = 8
    13 = 9
    14 = 17
    15 = 4
    16 = 10
    17 = 5
}

# Snapshot the original values for every row referenced above before writing
# anything, since several target rows borrow from rows that are themselves
# overwritten later in the loop.
$snapshot = @{}
foreach ($r in 2..17) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $sourceVals = $snapshot[$sourceRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $sourceVals[$c]
    }
}
